$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell that keeps its original formatting (style index 3),
# used to restore formatting on percentage cells after text assignment.
$formatDonor = $ws.Range("H2")

$ws.Range("E2").Value = "2026-02-06 09:48:04"
$ws.Range("K2").Value = "0.7 MJ/m2"
$ws.Range("M2").Value = "-0.3 °C 9:29 TU"
$ws.Range("E3").Value = "2026-02-06 09:48:07"
$ws.Range("K3").Value = "0.7 MJ/m2"
$ws.Range("M3").Value = "-1.1 °C 9:01 TU"
$ws.Range("O3").Value = "-2.6 °C"
$ws.Range("E4").Value = "2026-02-06 09:48:10"
$ws.Range("J4").Value = "995.1 hPa"
$ws.Range("K4").Value = "1.4 MJ/m2"
$ws.Range("L4").Value = "39.2 km/h - 299º 9:02 TU"
$ws.Range("O4").Value = "11.7 °C"
$ws.Range("E5").Value = "2026-02-06 09:48:13"
$ws.Range("J5").Value = "995.4 hPa"
$ws.Range("K5").Value = "1.4 MJ/m2"
$ws.Range("M5").Value = "11.6 °C 9:29 TU"
$ws.Range("O5").Value = "8.0 °C"
$ws.Range("E6").Value = "2026-02-06 09:48:15"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "51%"
$formatDonor.Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("J6").Value = "996.6 hPa"
$ws.Range("K6").Value = "1.7 MJ/m2"
$ws.Range("M6").Value = "15.7 °C 9:29 TU"
$ws.Range("O6").Value = "14.2 °C"
$ws.Range("E7").Value = "2026-02-06 09:48:18"
$ws.Range("J7").Value = "996.3 hPa"
$ws.Range("K7").Value = "1.7 MJ/m2"
$ws.Range("M7").Value = "11.0 °C 9:29 TU"
$ws.Range("E8").Value = "2026-02-06 09:48:21"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "91%"
$formatDonor.Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("K8").Value = "1.6 MJ/m2"
$ws.Range("L8").Value = "10.4 km/h - 196º 9:26 TU"
$ws.Range("M8").Value = "12.9 °C 9:21 TU"
$ws.Range("O8").Value = "6.2 °C"
$ws.Range("E9").Value = "2026-02-06 09:48:24"
$ws.Range("E10").Value = "2026-02-06 09:48:26"
$ws.Range("M10").Value = "9.4 °C 9:29 TU"
$ws.Range("O10").Value = "5.1 °C"
$ws.Range("E11").Value = "2026-02-06 09:48:29"
$ws.Range("J11").Value = "997.7 hPa"
$ws.Range("K11").Value = "0.5 MJ/m2"
$ws.Range("E12").Value = "2026-02-06 09:48:32"
$ws.Range("K12").Value = "1.4 MJ/m2"
$ws.Range("M12").Value = "16.2 °C 9:29 TU"
$ws.Range("O12").Value = "12.3 °C"
$ws.Range("E13").Value = "2026-02-06 09:48:34"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "90%"
$formatDonor.Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("M13").Value = "14.1 °C 9:20 TU"
$ws.Range("O13").Value = "6.3 °C"
$ws.Range("E14").Value = "2026-02-06 09:48:37"
$ws.Range("K14").Value = "0.3 MJ/m2"
$ws.Range("E15").Value = "2026-02-06 09:48:39"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "88%"
$formatDonor.Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("I15").Value = "0.2 mm"
$ws.Range("J15").Value = "995.6 hPa"
$ws.Range("K15").Value = "1.6 MJ/m2"
$ws.Range("M15").Value = "14.1 °C 9:29 TU"
$ws.Range("O15").Value = "6.1 °C"
$ws.Range("E16").Value = "2026-02-06 09:48:42"
$ws.Range("K16").Value = "1.3 MJ/m2"
$ws.Range("M16").Value = "5.8 °C 9:28 TU"
$ws.Range("E17").Value = "2026-02-06 09:48:45"
$ws.Range("J17").Value = "998.3 hPa"
$ws.Range("K17").Value = "1.4 MJ/m2"
$ws.Range("M17").Value = "5.8 °C 9:27 TU"
$ws.Range("O17").Value = "2.9 °C"
$ws.Range("E18").Value = "2026-02-06 09:48:47"
$ws.Range("K18").Value = "0.8 MJ/m2"
$ws.Range("E19").Value = "2026-02-06 09:48:51"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "91%"
$formatDonor.Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("J19").Value = "998.5 hPa"
$ws.Range("K19").Value = "1.4 MJ/m2"
$ws.Range("M19").Value = "11.3 °C 9:29 TU"
$ws.Range("O19").Value = "7.1 °C"
$ws.Range("E20").Value = "2026-02-06 09:48:54"
$ws.Range("K20").Value = "2.1 MJ/m2"
$ws.Range("E21").Value = "2026-02-06 09:48:56"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "87%"
$formatDonor.Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("J21").Value = "996.4 hPa"
$ws.Range("K21").Value = "1.7 MJ/m2"
$ws.Range("M21").Value = "12.0 °C 9:29 TU"
$ws.Range("O21").Value = "4.6 °C"
$ws.Range("E22").Value = "2026-02-06 09:48:59"
$ws.Range("K22").Value = "1.7 MJ/m2"
$ws.Range("O22").Value = "7.4 °C"
$ws.Range("E23").Value = "2026-02-06 09:49:01"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "93%"
$formatDonor.Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("J23").Value = "995.6 hPa"
$ws.Range("K23").Value = "1.5 MJ/m2"
$ws.Range("M23").Value = "10.4 °C 9:24 TU"
$ws.Range("O23").Value = "7.1 °C"
$ws.Range("E24").Value = "2026-02-06 09:49:04"
$ws.Range("J24").Value = "994.6 hPa"
$ws.Range("K24").Value = "1.4 MJ/m2"
$ws.Range("O24").Value = "12.3 °C"
$ws.Range("E25").Value = "2026-02-06 09:49:06"
$ws.Range("J25").Value = "997.6 hPa"
$ws.Range("K25").Value = "0.9 MJ/m2"
$ws.Range("M25").Value = "5.0 °C 9:25 TU"
$ws.Range("O25").Value = "2.1 °C"
$ws.Range("E26").Value = "2026-02-06 09:49:09"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "79%"
$formatDonor.Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("K26").Value = "0.5 MJ/m2"
$ws.Range("E27").Value = "2026-02-06 09:49:12"
$ws.Range("J27").Value = "995.4 hPa"
$ws.Range("K27").Value = "1.2 MJ/m2"
$ws.Range("O27").Value = "7.2 °C"
$ws.Range("E28").Value = "2026-02-06 09:49:15"
$ws.Range("J28").Value = "998.9 hPa"
$ws.Range("O28").Value = "1.5 °C"
$ws.Range("E29").Value = "2026-02-06 09:49:18"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "70%"
$formatDonor.Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("K29").Value = "1.5 MJ/m2"
$ws.Range("M29").Value = "14.8 °C 9:26 TU"
$ws.Range("O29").Value = "10.0 °C"
$ws.Range("E30").Value = "2026-02-06 09:49:20"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "70%"
$formatDonor.Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("K30").Value = "1.8 MJ/m2"
$ws.Range("E31").Value = "2026-02-06 09:49:21"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "99%"
$formatDonor.Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("J31").Value = "998.2 hPa"
$ws.Range("M31").Value = "6.3 °C 9:29 TU"
$ws.Range("E32").Value = "2026-02-06 09:49:23"
$ws.Range("J32").Value = "996.9 hPa"
$ws.Range("K32").Value = "1.5 MJ/m2"
$ws.Range("M32").Value = "16.9 °C 9:01 TU"
$ws.Range("O32").Value = "14.8 °C"
$ws.Range("E33").Value = "2026-02-06 09:49:24"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "96%"
$formatDonor.Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("M33").Value = "12.5 °C 9:23 TU"
$ws.Range("O33").Value = "6.6 °C"
$ws.Range("E34").Value = "2026-02-06 09:49:25"
$ws.Range("K34").Value = "1.5 MJ/m2"
$ws.Range("O34").Value = "6.0 °C"
$ws.Range("E35").Value = "2026-02-06 09:49:26"
$ws.Range("K35").Value = "0.6 MJ/m2"
$ws.Range("M35").Value = "-1.9 °C 9:29 TU"
$ws.Range("E36").Value = "2026-02-06 09:49:27"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "70%"
$formatDonor.Copy()
$ws.Range("H36").PasteSpecial(-4122)
$ws.Range("J36").Value = "998.3 hPa"
$ws.Range("K36").Value = "1.4 MJ/m2"
$ws.Range("M36").Value = "15.0 °C 9:18 TU"
$ws.Range("O36").Value = "10.7 °C"

$excel.CutCopyMode = 0

